$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasting Notes")

# Add a new row (row 15) with the new tasting note entry.
$ws.Range("A15").Value = 44171
$ws.Range("B15").Formula = '=A15-$A$6'
$ws.Range("C15").Value = 2.75
$ws.Range("D15").Value = "Cooled and served @ 12 C. Opening the swing-top gave a small puff. Poured clear with a foam. Moderate carbonation. A light sour grassy taste, less boozy and a dry mouthfeel."
